$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = 2
$ws.Range("C8").Value = 2

$ws.Range("H17").Select()
